$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "22×63=" "14×11="
Replace-Text "19×31=" "95×98="
Replace-Text "90×65=" "88×44="
Replace-Text "52×45=" "96×67="
Replace-Text "82×22=" "90×38="
Replace-Text "95×99=" "95×96="
Replace-Text "83×50=" "17×21="
Replace-Text "98×25=" "53×38="
Replace-Text "78×46=" "59×72="
Replace-Text "72×71=" "99×72="
Replace-Text "78×90=" "27×30="
Replace-Text "47×81=" "42×72="
Replace-Text "51×53=" "84×12="
Replace-Text "91×47=" "79×24="
Replace-Text "57×31=" "27×38="
Replace-Text "65×20=" "92×87="
Replace-Text "61×51=" "54×16="
Replace-Text "32×88=" "48×59="
Replace-Text "44×88=" "41×93="
Replace-Text "72×88=" "61×15="
Replace-Text "73×75=" "89×87="
Replace-Text "22×79=" "67×60="
Replace-Text "34×58=" "14×80="
Replace-Text "43×70=" "12×90="
Replace-Text "32×23=" "40×33="
